$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-CellText "D2" '245.70'
Set-CellText "D3" '22.12'
Set-CellText "D4" '5.368'
Set-CellText "D6" '3.396'
Set-CellText "D7" '6.362'
Set-CellText "D8" '0.8131'
Set-CellText "D9" '1.028'
Set-CellText "B10" 'One'
Set-CellText "C10" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-CellText "D10" '0.01121'
Set-CellText "E10" '9OneONEBestin24h'
Set-CellText "B11" 'WazirX'
Set-CellText "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-CellText "D11" '0.1422'
Set-CellText "E11" '10WazirXWRX'
Set-CellText "B12" 'LiechtensteinCryptoassetsExchange'
Set-CellText "C12" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-CellText "D12" '0.04127'
Set-CellText "E12" '11LiechtensteinCryptoassetsExchangeLCX'
Set-CellText "B13" 'MandalaExchangeToken'
Set-CellText "C13" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-CellText "D13" '0.07390'
Set-CellText "E13" '12MandalaExchangeTokenMDX'
Set-CellText "B14" 'BitrueCoin'
Set-CellText "C14" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-CellText "D14" '0.02986'
Set-CellText "E14" '13BitrueCoinBTR'
Set-CellText "B15" 'MCDex'
Set-CellText "C15" 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-CellText "D15" '4.148'
Set-CellText "E15" '14MCDexMCB'
Set-CellText "B16" 'BitMartToken'
Set-CellText "C16" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-CellText "D16" '0.09399'
Set-CellText "E16" '15BitMartTokenBMX'
Set-CellText "B17" 'BitForexToken'
Set-CellText "C17" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-CellText "D17" '0.001590'
Set-CellText "E17" '16BitForexTokenBF'
Set-CellText "B18" 'CoinExToken'
Set-CellText "C18" 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-CellText "D18" '0.04818'
Set-CellText "E18" '17CoinExTokenCET'
Set-CellText "D19" '0.006084'
Set-CellText "D20" '0.004085'
Set-CellText "D21" '0.0009839'
Set-CellText "D22" '0.0001501'
Set-CellText "D23" '3.716'
Set-CellText "D24" '2.229'
Set-CellText "D25" '0.3240'
Set-CellText "D27" '0.0002484'
Set-CellText "D40" '0.03876'
Set-CellText "D41" '0.006403'
Set-CellText "D43" '0.002601'
Set-CellText "D44" '0.005074'
Set-CellText "D45" '0.00005629'
Set-CellText "D47" '0.9204'
Set-CellText "E47" '46CoinbaseStockTokenCOIN'
Set-CellText "D49" '0.00002101'
